# [fix] create computer fixed
# Re-creates the "Entidade C" sub-level column and the "Tipo de computador"
# column, refreshes a handful of Entidade/Linha-type values, turns the
# e-mail in B11 into a mailto hyperlink with an updated address, and moves
# the active selection the way the author left the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the two new columns.
#    - a new "Entidade C" sub-column right after the existing Entidade C (old col G "Linha" shifts right)
#    - a new "Tipo de computador" column right before "Serial Number do Notebook" (old col N shifts right)
$ws.Columns("G").Insert()
$ws.Columns("O").Insert()

# 2) Header row touch-up for the "Entidade C" sub-column (reuses the existing
#    shared string). The "Tipo de computador" header is filled in further
#    down, in its original authoring order.
$ws.Range("G1").Value2 = "Entidade C"

# 3) Data updates - written in the same order the values were first typed
#    so the shared-string table grows the way it did originally.
$ws.Range("B11").Value2 = "Cr@email.com"

$ws.Range("F2").Value2 = "Volkswagen"
$ws.Range("F3").Value2 = "Volkswagen"

$ws.Range("F4").Value2 = "Saint Gobain"
$ws.Range("F5").Value2 = "Saint Gobain"
$ws.Range("F6").Value2 = "Saint Gobain"

$ws.Range("G4").Value2 = "Promotor"
$ws.Range("G6").Value2 = "Promotor"

$ws.Range("G5").Value2 = "Tecnico"

$ws.Range("F11").Value2 = "Solution Center"

$ws.Range("F7").Value2 = "Gomes da Costa"

$ws.Range("F8").Value2 = "Samsung"
$ws.Range("F9").Value2 = "Samsung"

$ws.Range("O1").Value2 = "Tipo de computador"

# "Tipo de computador" values (plain numbers)
$ws.Range("O2").Value2 = 1
$ws.Range("O3").Value2 = 1
$ws.Range("O4").Value2 = 1
$ws.Range("O7").Value2 = 2
$ws.Range("O8").Value2 = 2
$ws.Range("O10").Value2 = 1

# 4) Turn Camila's e-mail into a hyperlink.
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:Cr@email.com")
try { $wb.Styles("Hyperlink").NameLocal = "Hiperlink" } catch {}
try { $wb.Styles("Hyperlink").Name = "Hiperlink" } catch {}

# 5) Column widths for the two freshly inserted columns / reflowed neighbours.
$ws.Columns("G").ColumnWidth = 9.59
$ws.Columns("O").ColumnWidth = 26.17

# 6) View state: drop the old frozen top-left cell, move the active cell.
$ws.Range("O16").Select()

Write-Output "done"
